$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.143.56'
$ws.Range('D3').Value = '1.863.79'
$ws.Range('E3').Value = '  -0.84%  '
$ws.Range('D4').Value = '0.9997'
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '0.7091'
$ws.Range('E5').Value = '  -1.07%  '
$ws.Range('D6').Value = '241.41'
$ws.Range('E6').Value = '  -0.20%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('D8').Value = '0.3092'
$ws.Range('E8').Value = '  -1.12%  '
$ws.Range('E9').Value = '  -3.45%  '
$ws.Range('D10').Value = '24.61'
$ws.Range('E10').Value = '  -2.94%  '
$ws.Range('D11').Value = '0.08358'
$ws.Range('E11').Value = '  +1.02%  '
$ws.Range('D12').Value = '1.871.94'
$ws.Range('E12').Value = '  -2.61%  '
$ws.Range('D13').Value = '5.210'
$ws.Range('E13').Value = '  -1.68%  '
$ws.Range('D14').Value = '0.7073'
$ws.Range('E14').Value = '  -3.47%  '
$ws.Range('D15').Value = '91.14'
$ws.Range('D16').Value = '29.171.90'
$ws.Range('E16').Value = '  -1.14%  '
$ws.Range('D17').Value = '5.902'
$ws.Range('E17').Value = '  -0.94%  '
$ws.Range('D18').Value = '242.42'
$ws.Range('E18').Value = '  -2.21%  '
$ws.Range('D19').Value = '0.000007796'
$ws.Range('E19').Value = '  -0.95%  '
$ws.Range('D20').Value = '2.112.56'
$ws.Range('E20').Value = '  -3.33%  '
$ws.Range('E21').Value = '  -2.31%  '
$ws.Range('D22').Value = '0.9995'
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('D23').Value = '7.860'
$ws.Range('E23').Value = '  -2.05%  '
$ws.Range('D24').Value = '0.9999'
$ws.Range('E24').Value = '  +0.11%  '
$ws.Range('E25').Value = '  -1.46%  '
$ws.Range('D26').Value = '164.13'
$ws.Range('E26').Value = '  +0.17%  '
$ws.Range('D27').Value = '8.944'
$ws.Range('E27').Value = '  -1.36%  '
$ws.Range('D28').Value = '18.42'
$ws.Range('E28').Value = '  +0.25%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').Value = '1.499'
$ws.Range('E29').Value = '  -0.26%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').Value = '1.322'
$ws.Range('E30').Value = '  -2.93%  '
$ws.Range('D31').Value = '4.389'
$ws.Range('E31').Value = '  -0.12%  '
$ws.Range('D32').Value = '4.236'
$ws.Range('E32').Value = '  +2.70%  '
$ws.Range('D33').Value = '0.05141'
$ws.Range('E33').Value = '  -2.82%  '
$ws.Range('D34').Value = '0.7927'
$ws.Range('E34').Value = '  +8.65%  '
$ws.Range('D35').Value = '1.908'
$ws.Range('E35').Value = '  -2.57%  '
$ws.Range('D36').Value = '1.162'
$ws.Range('E36').Value = '  -3.35%  '
$ws.Range('D37').Value = '2.687'
$ws.Range('E37').Value = '  +0.48%  '
$ws.Range('D38').Value = '0.01845'
$ws.Range('E38').Value = '  -1.52%  '
$ws.Range('D39').Value = '2.699'
$ws.Range('E39').Value = '  -1.29%  '
$ws.Range('D40').Value = '1.164.23'
$ws.Range('E40').Value = '  -5.11%  '
$ws.Range('D41').Value = '6.230'
$ws.Range('E41').Value = '  +0.29%  '
$ws.Range('D42').Value = '0.8891'
$ws.Range('E42').Value = '  -2.48%  '
$ws.Range('D43').Value = '72.83'
$ws.Range('E43').Value = '  -2.96%  '
$ws.Range('D44').Value = '0.9996'
$ws.Range('E44').Value = '  +0.06%  '
$ws.Range('D45').Value = '102.62'
$ws.Range('E45').Value = '  -0.19%  '
$ws.Range('D46').Value = '2.009.10'
$ws.Range('E46').Value = '  -2.00%  '
$ws.Range('D47').Value = '0.5180'
$ws.Range('E47').Value = '  -1.58%  '
$ws.Range('E48').Value = '  -0.43%  '
$ws.Range('D49').Value = '9.277'
$ws.Range('E49').Value = '  -0.87%  '
$ws.Range('D50').Value = '1.002'
$ws.Range('E50').Value = '  -0.37%  '
$ws.Range('E51').Value = '  -1.76%  '
